# Add tp 26 and 27 from OR: extend the CRM accuracy table with a new row (80)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 80 data (mirrors the pattern of the preceding rows 77:79)
$ws.Range("A80").Value = 43815
$ws.Range("B80").Value = 2202.4457071062802
$ws.Range("C80").Value = 2207.0300000000002
$ws.Range("E80").Value = 169
$ws.Range("F80").Value = "New CRM opened 12/11/2022"

# Match the date formatting used by the rest of column A (copy from A79)
$ws.Range("A79").Copy()
$ws.Range("A80").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# Continue the "% off" formula series into the new row
$ws.Range("D80").Formula = "=100*(B80-C80)/C80"

# Update the view to reflect the new scroll position / selection
$excel.ActiveWindow.ScrollRow = 67
$ws.Range("A81").Select()
